# dadosCovid19.xlsx - "Add files via upload"
#
# Appends 15 new daily-tracking rows (2022/11/02 .. 2022/11/16) to the bottom
# of Sheet1, right after the existing last row (948, dated 2022/11/01).
# Columns follow the sheet's existing sparse layout: A,B,C,D,F,G,H,I,J,K,L,M,
# O,P,Q,T are populated; E,N,R,S are left blank, exactly like the surrounding
# rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# One row per hashtable. Column A is omitted from the "value dictionary"
# loop below and handled specially so the YYYY/MM/DD text isn't
# auto-recognised by Excel as a date serial.
$newRows = @(
    @{ Row=949; A="2022/11/02"; B=948; C=39448; D=691; F=0.01751673089; G=38685; H=27040; I=66488; J=1; K=1; L=1; M=0; O=0; P=0; Q=2; T=136 },
    @{ Row=950; A="2022/11/03"; B=949; C=39448; D=691; F=0.01751673089; G=38685; H=27040; I=66488; J=1; K=1; L=1; M=0; O=0; P=0; Q=2; T=136 },
    @{ Row=951; A="2022/11/04"; B=950; C=39448; D=691; F=0.01751673089; G=38685; H=27040; I=66488; J=1; K=1; L=1; M=0; O=0; P=0; Q=2; T=136 },
    @{ Row=952; A="2022/11/05"; B=951; C=39448; D=691; F=0.01751673089; G=38685; H=27040; I=66488; J=1; K=1; L=1; M=0; O=0; P=0; Q=2; T=136 },
    @{ Row=953; A="2022/11/06"; B=952; C=39448; D=691; F=0.01751673089; G=38685; H=27040; I=66488; J=1; K=1; L=1; M=0; O=0; P=0; Q=2; T=137 },
    @{ Row=954; A="2022/11/07"; B=953; C=39448; D=691; F=0.01751673089; G=38685; H=27040; I=66488; J=1; K=1; L=1; M=0; O=0; P=0; Q=2; T=137 },
    @{ Row=955; A="2022/11/08"; B=954; C=39448; D=691; F=0.01751673089; G=38685; H=27040; I=66488; J=1; K=1; L=1; M=0; O=0; P=0; Q=2; T=137 },
    @{ Row=956; A="2022/11/09"; B=955; C=39448; D=691; F=0.01751673089; G=38685; H=27040; I=66488; J=1; K=1; L=1; M=0; O=0; P=0; Q=2; T=137 },
    @{ Row=957; A="2022/11/10"; B=956; C=39448; D=691; F=0.01751673089; G=38685; H=27040; I=66488; J=1; K=1; L=1; M=0; O=0; P=0; Q=2; T=137 },
    @{ Row=958; A="2022/11/11"; B=957; C=39448; D=691; F=0.01751673089; G=38685; H=27040; I=66488; J=1; K=1; L=1; M=0; O=0; P=0; Q=2; T=137 },
    @{ Row=959; A="2022/11/12"; B=958; C=39448; D=691; F=0.01751673089; G=38685; H=27040; I=66488; J=1; K=1; L=1; M=0; O=0; P=0; Q=2; T=137 },
    @{ Row=960; A="2022/11/13"; B=959; C=39448; D=691; F=0.01751673089; G=38685; H=27040; I=66488; J=1; K=1; L=1; M=0; O=0; P=0; Q=2; T=138 },
    @{ Row=961; A="2022/11/14"; B=960; C=39448; D=691; F=0.01751673089; G=38685; H=27040; I=66488; J=1; K=1; L=1; M=0; O=0; P=0; Q=2; T=138 },
    @{ Row=962; A="2022/11/15"; B=961; C=39448; D=691; F=0.01751673089; G=38685; H=27040; I=66488; J=1; K=1; L=1; M=0; O=0; P=0; Q=2; T=138 },
    @{ Row=963; A="2022/11/16"; B=962; C=39486; D=691; F=0.01749987337; G=38709; H=27040; I=66526; J=38;          K=0; L=38; O=0; P=1; Q=1; T=138 }
)

# Map the data-dictionary keys (minus "Row"/"A") to their 1-based column
# numbers, in the same left-to-right order the original rows use.
$colNumber = @{ B=2; C=3; D=4; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13; O=15; P=16; Q=17; T=20 }
$colOrder  = @("B","C","D","F","G","H","I","J","K","L","M","O","P","Q","T")

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    # Column A holds a "YYYY/MM/DD" label that must stay plain text (like
    # every other date label already in the sheet) instead of being
    # reinterpreted as a date serial number. Temporarily force a text
    # number format while the literal is entered, then drop back to the
    # sheet's default (unstyled) formatting.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value2 = $rowData.A
    $ws.Cells.Item($r, 1).ClearFormats()

    foreach ($col in $colOrder) {
        if ($rowData.ContainsKey($col)) {
            $ws.Cells.Item($r, $colNumber[$col]).Value2 = $rowData[$col]
        }
    }
}
